$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45185
    }
}
